$d = $word.ActiveDocument

# 1. Programa section - split run text into numbered lines using manual line breaks (vertical tab => <w:br/>)
#    Two consecutive breaks are used between the major numbered sections (1..5) to reproduce the blank line.
$r1 = $d.Content
$found1 = $r1.Find.Execute("1. Conceitos ligados ao escoamento de fluídos e equações fundamentais1.1. Características e definições dos escoamentos;1.2. Conceitos de sistema e volume de controle;1.3. Equação da conservação da massa;1.4. Equação da conservação da energia;1.5. Equação da conservação da quantidade de movimento;1.6. Introdução à análise diferencial do movimento de fluidos.2. Escoamento incompressível de fluidos não viscosos2.1. Equação de Euler;2.2. Equação de Bernoulli;2.3. Aplicações da equação de Bernoulli.3. Escoamento viscoso incompressível3.1. Atrito e perda de carga;3.2. Avaliação das perdas de carga: regime laminar e turbulento;3.3. Equações de Hagen - Poiseuille e Darcy – Weisbach3.4. Diagrama de Moody e Moody –Rouse;3.5. Método dos comprimentos equivalentes.3.6. Presença de máquina no escoamento (bomba e turbina), Potência e rendimento;3.7. Medidores de vazão.4. Transferência de Calor4.1. Definição de Calor.4.2. Mecanismo da Condução.4.3. Mecanismo da Convecção.4.4. Associação de Mecanismos.5. Transferência de Massa5.1. Difusão e convecção mássica;5.2. 1ª lei de Fick;5.3. Concentrações mássica e molar;5.4. Frações mássica e molar;5.5. Velocidades médias mássica e molar;5.6. Fluxos difusivo mássico, difusivo molar, convectivo mássico e convectivo molar;5.7. Fluxo mássico total e fluxo molar total.", $false)
if (-not $found1) { throw "Programa text not found" }
$r1.Text = "1. Conceitos ligados ao escoamento de fluídos e equações fundamentais`v1.1. Características e definições dos escoamentos;`v1.2. Conceitos de sistema e volume de controle;`v1.3. Equação da conservação da massa;`v1.4. Equação da conservação da energia;`v1.5. Equação da conservação da quantidade de movimento;`v1.6. Introdução à análise diferencial do movimento de fluidos.`v`v2. Escoamento incompressível de fluidos não viscosos`v2.1. Equação de Euler;`v2.2. Equação de Bernoulli;`v2.3. Aplicações da equação de Bernoulli.`v`v3. Escoamento viscoso incompressível`v3.1. Atrito e perda de carga;`v3.2. Avaliação das perdas de carga: regime laminar e turbulento;`v3.3. Equações de Hagen - Poiseuille e Darcy – Weisbach`v3.4. Diagrama de Moody e Moody –Rouse;`v3.5. Método dos comprimentos equivalentes.`v3.6. Presença de máquina no escoamento (bomba e turbina), Potência e rendimento;`v3.7. Medidores de vazão.`v`v4. Transferência de Calor`v4.1. Definição de Calor.`v4.2. Mecanismo da Condução.`v4.3. Mecanismo da Convecção.`v4.4. Associação de Mecanismos.`v`v5. Transferência de Massa`v5.1. Difusão e convecção mássica;`v5.2. 1ª lei de Fick;`v5.3. Concentrações mássica e molar;`v5.4. Frações mássica e molar;`v5.5. Velocidades médias mássica e molar;`v5.6. Fluxos difusivo mássico, difusivo molar, convectivo mássico e convectivo molar;`v5.7. Fluxo mássico total e fluxo molar total."

# 2. Avaliacao / Criterio run: split "Nota de duas provas (P1 e P2)" from "Formula: M1 = ..."
$r2 = $d.Content
$found2 = $r2.Find.Execute("Nota de duas provas (P1 e P2)Fórmula: M1 = (P1 + 2 x P2)/3..", $false)
if (-not $found2) { throw "Avaliacao criterio text not found" }
$r2.Text = "Nota de duas provas (P1 e P2)`vFórmula: M1 = (P1 + 2 x P2)/3.."

# 3. Avaliacao / Norma de recuperacao run: split "Aplicacao..." from "NR (nota..."
$r3 = $d.Content
$found3 = $r3.Find.Execute("Aplicação de uma prova envolvendo o assunto de todo semestre.NR (nota da recuperação) = (M1 + NR)/2.", $false)
if (-not $found3) { throw "Avaliacao norma text not found" }
$r3.Text = "Aplicação de uma prova envolvendo o assunto de todo semestre.`vNR (nota da recuperação) = (M1 + NR)/2."

# 4. Bibliografia - split into 6 numbered reference lines
$r4 = $d.Content
$found4 = $r4.Find.Execute("1. FOX, R.W., MCDONALD, A.T., “Introdução à Mecânica dos Fluidos”, Ed. Guanabara Koogan.2. STREETER, V.L., WYLE,E.B., “Mecânica dos Fluidos”, Ed. Mc Graw Hill.3. OZISIK,M.N., “Transferência de Calor.”, Ed. Guanabara Koogan.4. INCROPERA, F.P.W., “Fundamentos de Transferência de Calor e Massa”, Ed. Guanabara Koogan.5. MUNSON, B.R.; YOUNG, D.F.; OKIISHI, T.H. Fundamentos da Mecânica dos Fluidos. Editora Edgard Blucher6 - GIORGETI, M. (2012) Fundamentos de Fenômenos de Transporte. Editora Campus", $false)
if (-not $found4) { throw "Bibliografia text not found" }
$r4.Text = "1. FOX, R.W., MCDONALD, A.T., “Introdução à Mecânica dos Fluidos”, Ed. Guanabara Koogan.`v2. STREETER, V.L., WYLE,E.B., “Mecânica dos Fluidos”, Ed. Mc Graw Hill.`v3. OZISIK,M.N., “Transferência de Calor.”, Ed. Guanabara Koogan.`v4. INCROPERA, F.P.W., “Fundamentos de Transferência de Calor e Massa”, Ed. Guanabara Koogan.`v5. MUNSON, B.R.; YOUNG, D.F.; OKIISHI, T.H. Fundamentos da Mecânica dos Fluidos. Editora Edgard Blucher`v6 - GIORGETI, M. (2012) Fundamentos de Fenômenos de Transporte. Editora Campus"

Write-Host "Done"
